$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = "test"
